# Apply the edit described by the diff: update the "Share" (column K) values
# for rows 10, 11, 13, 14, 15 from 1 to 0, and move the active selection to K10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0

$ws.Range("K10").Select()
